# Daily attendance processing - 2025-12-05 03:13:04
# Swap order of "System" and the reporter's email in column G ("Recorded By")
# for rows where the value is exactly "System, <email>" with a single email
# belonging to dnasr281@gmail.com or admin@admin.com.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value()

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
    }
}
